# Update "想去人数" (want-to-go count) values in F column on two sheets,
# mirroring the source-data refresh described in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 223
$wsExpo.Range("F15").Value = 29
$wsExpo.Range("F22").Value = 1630
$wsExpo.Range("F23").Value = 3811
$wsExpo.Range("F27").Value = 1137
$wsExpo.Range("F28").Value = 127
$wsExpo.Range("F29").Value = 2024

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F16").Value = 29
$wsAll.Range("F23").Value = 1630
$wsAll.Range("F24").Value = 3811
$wsAll.Range("F28").Value = 1137
$wsAll.Range("F29").Value = 127
$wsAll.Range("F30").Value = 2024
